{"js": "// Change nomenclature in \"verifikasi ka\" (KA administrative-completeness\n// table): the PIPPIB justification row used to read\n//   \"Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB\"\n// and is renamed to\n//   \"Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB\"\n// There are similarly worded rows elsewhere in the table (with spaces\n// around the \"/\" and without \"PIPPIB\"), so search on the full, unique\n// original sentence rather than just the \"persetujuan awal\" fragment.\n\nconst oldText =\n  \"Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB\";\nconst newText =\n  \"Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldText);\n}\n\n// Replace the matched range's text in place so the run keeps its\n// existing formatting (Tahoma, 8pt/sz16, noProof).\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change nomenclature in \"verifikasi ka\" (KA administrative-completeness\n# table): the PIPPIB justification row used to read\n#   \"Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB\"\n# and is renamed to\n#   \"Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB\"\n#\n# The table has two other, similarly worded rows (\"Justifikasi / bukti\n# kesesuaian lokasi ... tata ruang ...\" and \"Justifikasi / bukti\n# persetujuan awal rencana usaha dan/atau kegiatan\" with no \"PIPPIB\"),\n# so we search on the full, unique original sentence (no spaces around\n# the \"/\") instead of just the \"persetujuan awal\" fragment, to make\n# sure only the intended row is touched.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB\"\n$newText = \"Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
